$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete the entire column A, shifting all remaining columns (B:N) one to the left (-> A:M)
$ws.Range("A:A").EntireColumn.Delete()
